$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 ("time_taken"), matching the style of the other
# header cells (bold/bordered/centered) by copying format from E1.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate F2:F30 with the recorded time_taken timestamps (plain strings,
# same unstyled formatting as the other data cells in those rows).
$ws.Range("F2").Value = "2021-10-05 10:52:07.246784"
$ws.Range("F3").Value = "2021-10-05 10:52:07.246796"
$ws.Range("F4").Value = "2021-10-05 10:52:07.246800"
$ws.Range("F5").Value = "2021-10-05 10:52:07.246804"
$ws.Range("F6").Value = "2021-10-05 10:52:07.246807"
$ws.Range("F7").Value = "2021-10-05 10:52:07.246811"
$ws.Range("F8").Value = "2021-10-05 10:52:07.246814"
$ws.Range("F9").Value = "2021-10-05 10:52:07.246817"
$ws.Range("F10").Value = "2021-10-05 10:52:07.246821"
$ws.Range("F11").Value = "2021-10-05 10:52:07.246824"
$ws.Range("F12").Value = "2021-10-05 10:52:07.246827"
$ws.Range("F13").Value = "2021-10-05 10:52:07.246830"
$ws.Range("F14").Value = "2021-10-05 10:52:07.246833"
$ws.Range("F15").Value = "2021-10-05 10:52:07.246837"
$ws.Range("F16").Value = "2021-10-05 10:52:07.246839"
$ws.Range("F17").Value = "2021-10-05 10:52:07.246843"
$ws.Range("F18").Value = "2021-10-05 10:52:07.246846"
$ws.Range("F19").Value = "2021-10-05 10:52:07.246849"
$ws.Range("F20").Value = "2021-10-05 10:52:07.246852"
$ws.Range("F21").Value = "2021-10-05 10:52:07.246856"
$ws.Range("F22").Value = "2021-10-05 10:52:07.246859"
$ws.Range("F23").Value = "2021-10-05 10:52:07.246862"
$ws.Range("F24").Value = "2021-10-05 10:52:07.246865"
$ws.Range("F25").Value = "2021-10-05 10:52:07.246868"
$ws.Range("F26").Value = "2021-10-05 10:52:07.246872"
$ws.Range("F27").Value = "2021-10-05 10:52:07.246875"
$ws.Range("F28").Value = "2021-10-05 10:52:07.246878"
$ws.Range("F29").Value = "2021-10-05 10:52:07.246881"
$ws.Range("F30").Value = "2021-10-05 10:52:07.246884"
